# Capacity Supply Curve.xlsx — edit CSC value for onshore wind to represent
# higher siting constraints compared to solar.
#
# "onshore wind es" lives in row 7 of the "CSC-CSCSoCECBiaSY" sheet
# (Share of Cost Effective Capacity Built in a Single Year). Its yearly
# values (columns B:AE, years 2021-2050) move from 0.3 -> 0.2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CSC-CSCSoCECBiaSY")

$ws.Range("B7:AE7").Value = 0.2

# Leave the workbook focused/selected on the edited sheet and range, as the
# author did when they made (and saved after) this edit.
$ws.Activate()
$ws.Range("B7:AE7").Select()
